# Daily attendance processing - 2025-11-23 23:24:54
#
# The "Recorded By" column (G) lists who recorded a session as a
# comma-separated list of names/emails. Entries that end with the
# literal token "System" (and do not already start with "System")
# need to have that trailing "System" token swapped with the first
# token in the list, e.g.:
#   "dnasr281@gmail.com, System"                -> "System, dnasr281@gmail.com"
#   "system, backup@backdoor.com, System"       -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Text

    if ($current -eq $null -or $current -eq "") {
        continue
    }

    $parts = $current.Split(",")
    if ($parts.Length -gt 1) {
        $first = $parts[0].Trim()
        $lastIdx = $parts.Length - 1
        $lastPart = $parts[$lastIdx].Trim()

        if ($lastPart.Equals("System") -and -not $first.Equals("System")) {
            $newParts = New-Object 'object[]' $parts.Length
            $newParts[0] = $lastPart
            for ($i = 1; $i -lt $lastIdx; $i++) {
                $newParts[$i] = $parts[$i].Trim()
            }
            $newParts[$lastIdx] = $first

            $newValue = [string]::Join(", ", $newParts)
            $cell.Value = $newValue
        }
    }
}
